# Updated project to the state it was in at the end of project 3:
# populate the "Estimated Effort" (column C) values for the backlog items,
# and move the active selection to C36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Estimated Effort" values in column C ---
$ws.Range("C5").Value  = "Little"
$ws.Range("C6").Value  = "Moderate"
$ws.Range("C7").Value  = "Moderate"
$ws.Range("C8").Value  = "Moderate"
$ws.Range("C9").Value  = "Moderate"
$ws.Range("C10").Value = "Very Heavy"
$ws.Range("C11").Value = "Little"
$ws.Range("C12").Value = "Moderate"

$ws.Range("C15").Value = "Little"

$ws.Range("C17").Value = "Moderate"
$ws.Range("C18").Value = "Little"
$ws.Range("C19").Value = "Little"
$ws.Range("C21").Value = "Moderate"
$ws.Range("C22").Value = "Little"
$ws.Range("C23").Value = "Heavy"
$ws.Range("C24").Value = "Moderate"
$ws.Range("C25").Value = "Moderate"
$ws.Range("C27").Value = "Moderate"
$ws.Range("C28").Value = "Little"
$ws.Range("C29").Value = "Heavy"
$ws.Range("C30").Value = "Moderate"
$ws.Range("C31").Value = "Moderate"

$ws.Range("C34").Value = "Little"
$ws.Range("C35").Value = "Little"

# --- Move the active selection (as recorded when the author saved) ---
$ws.Range("C36").Select() | Out-Null
